$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4) mirroring the existing rows' layout.
$ws.Cells.Item(4, 1).Value = 42607.891712962963
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(4, 2).Value = -32
$ws.Cells.Item(4, 3).Value = 42
$ws.Cells.Item(4, 4).Value = 54
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 94
$ws.Cells.Item(4, 7).Value = 19078
$ws.Cells.Item(4, 8).Value = 14573
$ws.Cells.Item(4, 9).Value = 2277
$ws.Cells.Item(4, 10).Value = 226
$ws.Cells.Item(4, 11).Value = 287
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 18
$ws.Cells.Item(4, 14).Value = "Noun"
